$wb = $excel.ActiveWorkbook

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12501862
$ws.Range("J17").Value = 12501862
$ws.Range("L17").Value = 37505586
$ws.Range("N17").Value = -37505922

# ALC!row81
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 79789.5
$ws.Range("J81").Value = 79789.5
$ws.Range("L81").Value = 79789.5
$ws.Range("N81").Value = -81785.5

# ALC!row84
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 79789.5
$ws.Range("J84").Value = 79789.5
$ws.Range("L84").Value = 239368.5
$ws.Range("N84").Value = -249352.5

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1127.0834
$ws.Range("I129").Value = 792.5
$ws.Range("J129").Value = 2800
$ws.Range("K129").Value = 2377.5
$ws.Range("L129").Value = 8400
$ws.Range("M129").Value = 2622.5
$ws.Range("N129").Value = -18400

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3598.6365
$ws.Range("I137").Value = 2159
$ws.Range("K137").Value = 6477
$ws.Range("M137").Value = -3927

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4205.282
$ws.Range("I138").Value = 2094.6
$ws.Range("K138").Value = 6283.799999999999
$ws.Range("M138").Value = -1143.799999999999

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6368.8
$ws.Range("I141").Value = 5868.364
$ws.Range("K141").Value = 17605.092
$ws.Range("M141").Value = -12425.092

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3931.2058
$ws.Range("I110").Value = 3546.2307
$ws.Range("K110").Value = 3546.2307
$ws.Range("M110").Value = -1501.2307

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8061.0835
$ws.Range("I132").Value = 8173.1665
$ws.Range("J132").Value = 7949
$ws.Range("K132").Value = 24519.4995
$ws.Range("L132").Value = 23847
$ws.Range("M132").Value = -21989.4995
$ws.Range("N132").Value = -28907

# BSM!row2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69879
$ws.Range("J2").Value = 69879
$ws.Range("L2").Value = 69879
$ws.Range("N2").Value = -70105

# BSM!row76
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 16333.333
$ws.Range("J76").Value = 16333.333
$ws.Range("L76").Value = 16333.333
$ws.Range("N76").Value = -16963.333

# BSM!row79
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 16333.333
$ws.Range("J79").Value = 16333.333
$ws.Range("L79").Value = 16333.333
$ws.Range("N79").Value = -18517.333

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1651.8334
$ws.Range("J94").Value = 2650
$ws.Range("L94").Value = 2650
$ws.Range("N94").Value = -3552

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5118.727
$ws.Range("I134").Value = 5080.5454
$ws.Range("K134").Value = 15241.6362
$ws.Range("M134").Value = -12706.6362

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I122").Value = 2550.88
$ws.Range("J122").Value = 2058.4167
$ws.Range("K122").Value = 7652.64
$ws.Range("L122").Value = 6175.250100000001
$ws.Range("M122").Value = -5202.64
$ws.Range("N122").Value = -11075.2501

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9219.933999999999
$ws.Range("I134").Value = 7163.5454
$ws.Range("K134").Value = 21490.6362
$ws.Range("M134").Value = -18955.6362

# CUL!row23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 257
$ws.Range("J23").Value = 149.5
$ws.Range("L23").Value = 448.5
$ws.Range("N23").Value = -918.5

# CUL!row52
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1440.3334
$ws.Range("J52").Value = 1440.3334
$ws.Range("L52").Value = 4321.0002
$ws.Range("N52").Value = -4853.0002

# CUL!row60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 837.5
$ws.Range("I60").Value = 1033.3334
$ws.Range("J60").Value = 250
$ws.Range("K60").Value = 3100.0002
$ws.Range("L60").Value = 750
$ws.Range("N60").Value = -1252

# CUL!row124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 950
$ws.Range("I124").Value = 950
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 2850
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 2060
$ws.Range("N124").ClearContents()

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13160405
$ws.Range("I131").Value = 125001270
$ws.Range("J131").Value = 2655.8235
$ws.Range("K131").Value = 375003810
$ws.Range("L131").Value = 7967.470499999999
$ws.Range("M131").Value = -374998770
$ws.Range("N131").Value = -18047.4705

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3510.2666
$ws.Range("I137").Value = 1958.625
$ws.Range("J137").Value = 5283.5713
$ws.Range("K137").Value = 5875.875
$ws.Range("L137").Value = 15850.7139
$ws.Range("M137").Value = -775.875
$ws.Range("N137").Value = -26050.7139

# CUL!row140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1360.3684
$ws.Range("I140").Value = 1079.2354
$ws.Range("K140").Value = 3237.7062
$ws.Range("M140").Value = 1942.2938

# GSM!row99
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10690
$ws.Range("I99").Value = 10690
$ws.Range("K99").Value = 10690
$ws.Range("M99").Value = -8444

# GSM!row123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 54599.75
$ws.Range("J123").Value = 54599.75
$ws.Range("L123").Value = 54599.75
$ws.Range("N123").Value = -59499.75

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3144.65
$ws.Range("I132").Value = 3035.7273
$ws.Range("J132").Value = 3658.1428
$ws.Range("K132").Value = 9107.1819
$ws.Range("L132").Value = 10974.4284
$ws.Range("M132").Value = -6577.1819
$ws.Range("N132").Value = -16034.4284

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4534.143
$ws.Range("I7").Value = 4528.222
$ws.Range("K7").Value = 4528.222
$ws.Range("M7").Value = -4416.222

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2667.375
$ws.Range("I16").Value = 2667.375
$ws.Range("K16").Value = 2667.375
$ws.Range("M16").Value = -2497.375

# LTW!row103
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 24249.75
$ws.Range("J103").Value = 24249.75
$ws.Range("L103").Value = 24249.75
$ws.Range("N103").Value = -26593.75

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4534.143
$ws.Range("I126").Value = 4528.222
$ws.Range("K126").Value = 13584.666
$ws.Range("M126").Value = -11114.666

# WVR!row15
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7498.143
$ws.Range("J15").Value = 7499
$ws.Range("L15").Value = 7499
$ws.Range("N15").Value = -8075

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9997.333000000001
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9997.333000000001
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1407.8
$ws.Range("I81").Value = 1407.8
$ws.Range("K81").Value = 2815.6
$ws.Range("M81").Value = -1754.6

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1407.8
$ws.Range("I84").Value = 1407.8
$ws.Range("K84").Value = 14078
$ws.Range("M84").Value = -8774
